# Insert a new weekly data row before the current row 328 (new row = 329..355
# shift down by one). Fill the inserted row with the new week's values while
# keeping the remaining cells identical to the row previously at 328.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the full contents of the row currently at 328 so we can re-use the
# fields that are NOT changing for the newly inserted row.
$srcRow = 328
$A = $ws.Cells.Item($srcRow, 1).Value()
$B = $ws.Cells.Item($srcRow, 2).Value()
$C = $ws.Cells.Item($srcRow, 3).Value()
$E = $ws.Cells.Item($srcRow, 5).Value()
$F = $ws.Cells.Item($srcRow, 6).Value()
$G = $ws.Cells.Item($srcRow, 7).Value()
$H = $ws.Cells.Item($srcRow, 8).Value()
$I = $ws.Cells.Item($srcRow, 9).Value()
$N = $ws.Cells.Item($srcRow, 14).Value()
$O = $ws.Cells.Item($srcRow, 15).Value()
$Q = $ws.Cells.Item($srcRow, 17).Value()
$R = $ws.Cells.Item($srcRow, 18).Value()

# Push row 328 (and everything below it) down by one row.
$ws.Rows.Item(328).Insert()

# Populate the freshly inserted row 328 with the new weekly record.
$ws.Cells.Item(328, 1).Value = $A
$ws.Cells.Item(328, 2).Value = $B
$ws.Cells.Item(328, 3).Value = $C
$ws.Cells.Item(328, 4).Value = 44461
$ws.Cells.Item(328, 5).Value = $E
$ws.Cells.Item(328, 6).Value = $F
$ws.Cells.Item(328, 7).Value = $G
$ws.Cells.Item(328, 8).Value = $H
$ws.Cells.Item(328, 9).Value = $I
$ws.Cells.Item(328, 10).Value = 370
$ws.Cells.Item(328, 11).Value = 4500
$ws.Cells.Item(328, 12).Value = 5000
$ws.Cells.Item(328, 13).Value = 4703
$ws.Cells.Item(328, 14).Value = $N
$ws.Cells.Item(328, 15).Value = $O
$ws.Cells.Item(328, 16).Value = 470
$ws.Cells.Item(328, 17).Value = $Q
$ws.Cells.Item(328, 18).Value = $R
